# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "Datos actualizados..." timestamp cell (A1)
$ws.Range("A1").Value = "Datos actualizados a 22 de Abril de 2020 a las 15:52"

# Row 4 - Estados Unidos
$ws.Range("D4").Value = 83008
$ws.Range("E4").Value = 690957
$ws.Range("G4").Value = 38
$ws.Range("H4").Value = 45356

# Row 40 - Noruega
$ws.Range("B40").Value = 7275
$ws.Range("C40").Value = 34
$ws.Range("E40").Value = 7060
$ws.Range("F40").Value = 54
$ws.Range("G40").Value = 1
$ws.Range("H40").Value = 183

# Row 65 - Barein
$ws.Range("B65").Value = 2009
$ws.Range("C65").Value = 36
$ws.Range("D65").Value = 1026
$ws.Range("E65").Value = 976

# Row 67 - Islandia
$ws.Range("B67").Value = 1785
$ws.Range("C67").Value = 7
$ws.Range("D67").Value = 1462
$ws.Range("E67").Value = 313

# Row 158 - Uganda
$ws.Range("D158").Value = 45
$ws.Range("E158").Value = 16

# Row 166 - Macao
$ws.Range("D166").Value = 26
$ws.Range("E166").Value = 19

# Row 180 - Malaui
$ws.Range("E180").Value = 17
$ws.Range("G180").Value = 1
$ws.Range("H180").Value = 3
